$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Cell value writes (order chosen to control shared-string append order) ----
$ws.Range('K5').Value2 = 'V1.04.1'
$ws.Range('L5').Value2 = 0.01361
$ws.Range('N5').Value2 = 'V1.04.1'
$ws.Range('O5').Value2 = 0.01468
$ws.Range('K6').Value2 = 'V1.04.2'
$ws.Range('L6').Value2 = 0.01382
$ws.Range('N6').Value2 = 'V1.04.2'
$ws.Range('O6').Value2 = 0.01436
$ws.Range('F7').Value2 = 'V1.04.2'
$ws.Range('G7').Value2 = 25.76465
$ws.Range('H7').Value2 = 0.00269
$ws.Range('I7').Value2 = 0.91974
$ws.Range('K7').Value2 = 'V1.04.3'
$ws.Range('L7').Value2 = 0.01487
$ws.Range('N7').Value2 = 'V1.04.3'
$ws.Range('O7').Value2 = 0.01728
$ws.Range('F8').Value2 = 'V1.04.3'
$ws.Range('G8').Value2 = 18.13266
$ws.Range('H8').Value2 = 0.01539
$ws.Range('I8').Value2 = 0.73134
$ws.Range('K8').Value2 = 'V1.04.4'
$ws.Range('N8').Value2 = 'V1.04.4'
$ws.Range('F9').Value2 = 'V1.04.4'
$ws.Range('G9').Value2 = 16.71928
$ws.Range('H9').Value2 = 0.02133
$ws.Range('I9').Value2 = 0.74704
$ws.Range('A12').Value2 = 'V1.03.7'
$ws.Range('S3').Value2 = 'Mix 1'
$ws.Range('T3').Value2 = 15.36356
$ws.Range('U3').Value2 = 0.83669
$ws.Range('S4').Value2 = 'Mix 2'
$ws.Range('T4').Value2 = 15.54205
$ws.Range('U4').Value2 = 0.83874
$ws.Range('T1').Value2 = 'O-Haze Performance'
$ws.Range('W1').Value2 = 'RESIDE-OTS Performance'
$ws.Range('T2').Value2 = 'PSNR'
$ws.Range('U2').Value2 = 'SSIM'
$ws.Range('W2').Value2 = 'PSNR'
$ws.Range('X2').Value2 = 'SSIM'
$ws.Range('Q4').Value2 = 0.04015
